$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price column D, Volume(1h) column E)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.308.68'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.17%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.872.36'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.21%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("E5").Value = '  -0.83%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.07%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07798'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3107'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.27%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.11'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.26%  '

$ws.Range("E11").Value = '  +0.04%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.885.62'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.233'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.12%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7176'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.42%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.25'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.43%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008381'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.35%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.122'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.55%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.315.82'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.20%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '240.25'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.126.45'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.65%  '

$ws.Range("E21").Value = '  -0.17%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.743'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.82%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1593'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.65%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.57'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.96%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.025'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.14%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.46'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.55%  '

$ws.Range("E29").Value = '  -0.14%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.418'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.23%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.354'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.64%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.243'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.59%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05356'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.13%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.941'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.85%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7485'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.52%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.173'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.683'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01875'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.241.84'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.06%  '

$ws.Range("E40").Value = '  +0.42%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.534'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.89%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8925'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.32%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '109.85'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.93%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.28'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.66%  '

$ws.Range("E45").Value = '  +0.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000130'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.87%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.017.93'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.59%  '

$ws.Range("E48").Value = '  +0.00%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.797'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.06%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.427'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.26%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4337'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.71%  '
